$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.894.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.87%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.743.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.55%  "

# Row 4
$ws.Range("E4").Value = "  +0.24%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.61%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.42%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.12%  "

# Row 8
$ws.Range("E8").Value = "  +1.95%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.771.33"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.94%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.46%  "

# Row 11
$ws.Range("E11").Value = "  +5.99%  "

# Row 12
$ws.Range("E12").Value = "  +3.25%  "

# Row 13
$ws.Range("E13").Value = "  +2.59%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.236.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.76%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.83%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.053.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.19%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000155"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.52%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.764.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.92%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.36%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.45%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "363.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.88%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.50%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.00%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.534"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.74%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.91%  "

# Row 26
$ws.Range("E26").Value = "  +5.90%  "

# Row 27
$ws.Range("E27").Value = "  +5.86%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.27%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0913"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +12.21%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.33%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.94%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +20.29%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "175.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.35%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.10%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.73"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.03%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.39%  "

# Row 37
$ws.Range("E37").Value = "  +8.77%  "

# Row 38
$ws.Range("E38").Value = "  +10.68%  "

# Row 39
$ws.Range("E39").Value = "  +10.81%  "

# Row 40
$ws.Range("E40").Value = "  +5.32%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "342.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.15%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.22%  "

# Row 43
$ws.Range("E43").Value = "  +11.74%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.57%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.66%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0601"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.81%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.654"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.99%  "

# Row 48
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0261"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.25%  "

# Row 49
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.22%  "

# Row 50
$ws.Range("E50").Value = "  +2.23%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.994"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.29%  "

